$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The CLAVE column (A) for each group of rows under a given subject should
# all repeat the same key as the first row of that group. A previous (bad)
# edit had given each sub-row of a group its own unique -- and wrong --
# key pulled from spurious shared strings appended at the end of the
# sharedStrings table. This restores the correct/shared key per group.

$groups = @(
    @{ Key = "A3";  Rows = @("A4","A5","A6","A7","A8") },
    @{ Key = "A9";  Rows = @("A10","A11","A12","A13","A14") },
    @{ Key = "A15"; Rows = @("A16","A17","A18","A19","A20") },
    @{ Key = "A38"; Rows = @("A39","A40","A41","A42","A43") },
    @{ Key = "A44"; Rows = @("A45","A46","A47","A48","A49") },
    @{ Key = "A50"; Rows = @("A51","A52","A53","A54","A55") }
)

foreach ($group in $groups) {
    $keyValue = $ws.Range($group.Key).Value()
    foreach ($cellRef in $group.Rows) {
        $ws.Range($cellRef).Value = $keyValue
    }
}

# Match the final selection/scroll state left by the author after the fix.
$ws.Range("A50").Select() | Out-Null
